# Test if all specified projects has three Issues
#
# "browse_issues" sheet used to just be a flat list of project keys
# ("Projects to browse" / TOUCAN / JETI / COALA). Turn it into a small
# two-column table: project key + how many issues that project needs to
# have for the test to pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("browse_issues")

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Project"
$ws.Range("B1").Value = "Number of Issues required"
$ws.Range("A1:C1").Font.Bold = $true

# --- TOUCAN needs 3 issues -------------------------------------------
$ws.Range("A2").Value = "TOUCAN"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"

# --- JETI (no issue count specified, just keep the key as text) ------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "JETI"
$ws.Range("C3").NumberFormat = "@"

# --- COALA unchanged ---------------------------------------------------
$ws.Range("A4").Value = "COALA"

$ws.Range("D14").Select()
